# dissertation_figures.pptx edit:
#  1) Update the cached "datetimeFigureOut" footer/date field text from
#     6/3/2013 to 6/10/2013 on the slide master and every slide layout.
#  2) On slide 14, split two flowchart-box captions into two runs each:
#       "Generate Edge's Vertices"     -> "Create Edge's " + "Vertices"
#       "Generate edge in the graph"   -> "Add edge "      + "in the graph"

$p = $ppt.ActivePresentation

# --- 1) Date placeholder text -------------------------------------------------
$oldDate = "6/3/2013"
$newDate = "6/10/2013"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout (CustomLayouts hang off the design's slide master)
$design = $p.Designs.Item(1)
$layouts = $design.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2) Slide 14 caption edits ------------------------------------------------
$s14 = $p.Slides.Item(14)

$shpVertices = $s14.Shapes.Item("Flowchart: Process 8")
$trV = $shpVertices.TextFrame.TextRange
$trV.Text = "Create Edge" + [char]0x2019 + "s Vertices"
$firstPartV = $trV.Characters(1, 14)
$firstPartV.Text = "Create Edge" + [char]0x2019 + "s "

$shpAddEdge = $s14.Shapes.Item("Flowchart: Process 74")
$trA = $shpAddEdge.TextFrame.TextRange
$trA.Text = "Add edge in the graph"
$firstPartA = $trA.Characters(1, 9)
$firstPartA.Text = "Add edge "

Write-Output "done"
